# Update quarterly cash-flow database: drop oldest quarter (1399/06),
# shift remaining quarters left, and populate the newest quarter (1401/12).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 8
$ws.Range("D8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# Row 9
$ws.Range("D9").Value = "1400-10-29 (2)"
$ws.Range("E9").Value = "1401-02-07 (9)"
$ws.Range("F9").Value = "1401-04-26 (4)"
$ws.Range("G9").Value = "1401-08-28 (4)"
$ws.Range("H9").Value = "1401-10-29 (2)"
$ws.Range("I9").Value = "1402-02-13 (9)"
$ws.Range("J9").Value = "1401-04-26 (2)"
$ws.Range("K9").Value = "1401-08-28 (2)"
$ws.Range("L9").Value = "1401-10-29"
$ws.Range("M9").Value = "1402-02-13 (2)"

# Row 12
$ws.Range("D12").Value = 170634
$ws.Range("E12").Value = 146164
$ws.Range("F12").Value = 377750
$ws.Range("G12").Value = 321488
$ws.Range("H12").Value = 431215
$ws.Range("I12").Value = 449386
$ws.Range("J12").Value = 361049
$ws.Range("K12").Value = 408709
$ws.Range("L12").Value = 703939
$ws.Range("M12").Value = 366308

# Row 13
$ws.Range("D13").Value = -67106
$ws.Range("E13").Value = 20011
$ws.Range("F13").Value = -20000
$ws.Range("G13").Value = -75163
$ws.Range("H13").Value = 68364
$ws.Range("I13").Value = -78620
$ws.Range("J13").Value = -15000
$ws.Range("K13").Value = -29713
$ws.Range("L13").Value = -73807
$ws.Range("M13").Value = -44078

# Row 14
$ws.Range("D14").Value = 103528
$ws.Range("E14").Value = 166175
$ws.Range("F14").Value = 357750
$ws.Range("G14").Value = 246325
$ws.Range("H14").Value = 499579
$ws.Range("I14").Value = 370766
$ws.Range("J14").Value = 346049
$ws.Range("K14").Value = 378996
$ws.Range("L14").Value = 630132
$ws.Range("M14").Value = 322230

# Row 16
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0

# Row 17
$ws.Range("D17").Value = -5692
$ws.Range("E17").Value = -17346
$ws.Range("F17").Value = -3680
$ws.Range("G17").Value = -29764
$ws.Range("H17").Value = -20449
$ws.Range("I17").Value = -409312
$ws.Range("J17").Value = -368010
$ws.Range("K17").Value = -75878
$ws.Range("L17").Value = 121089
$ws.Range("M17").Value = -292809

# Row 18
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0

# Row 19
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0

# Row 20
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = -165
$ws.Range("G20").Value = -8649
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = -6224
$ws.Range("J20").Value = -31
$ws.Range("K20").Value = -7808
$ws.Range("L20").Value = -88
$ws.Range("M20").Value = -42

# Row 21
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0

# Row 22
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0

# Row 23
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0

# Row 25
$ws.Range("D25").Value = 56000
$ws.Range("E25").Value = 1454
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 110000
$ws.Range("J25").Value = 155000
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 161973
$ws.Range("M25").Value = 178266

# Row 26
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = -50000
$ws.Range("F26").Value = -100000
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = -150000
$ws.Range("I26").Value = -162656
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -118356
$ws.Range("M26").Value = -199441

# Row 27
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0

# Row 28
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0

# Row 29
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0

# Row 30
$ws.Range("D30").Value = 394
$ws.Range("E30").Value = 15126
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 36330
$ws.Range("J30").Value = 41265
$ws.Range("K30").Value = -33663
$ws.Range("L30").Value = 1
$ws.Range("M30").Value = 41528

# Row 31
$ws.Range("D31").Value = 18568
$ws.Range("E31").Value = 16741
$ws.Range("F31").Value = 20666
$ws.Range("G31").Value = 23277
$ws.Range("H31").Value = 28128
$ws.Range("I31").Value = 26202
$ws.Range("J31").Value = 25826
$ws.Range("K31").Value = 25203
$ws.Range("L31").Value = 131579
$ws.Range("M31").Value = -97842

# Row 32
$ws.Range("D32").Value = 69270
$ws.Range("E32").Value = -34025
$ws.Range("F32").Value = -83179
$ws.Range("G32").Value = -15136
$ws.Range("H32").Value = -142321
$ws.Range("I32").Value = -405660
$ws.Range("J32").Value = -145950
$ws.Range("K32").Value = -72146
$ws.Range("L32").Value = 296198
$ws.Range("M32").Value = -370340

# Row 33
$ws.Range("D33").Value = 172798
$ws.Range("E33").Value = 132150
$ws.Range("F33").Value = 274571
$ws.Range("G33").Value = 231189
$ws.Range("H33").Value = 357258
$ws.Range("I33").Value = -34894
$ws.Range("J33").Value = 200099
$ws.Range("K33").Value = 306850
$ws.Range("L33").Value = 926330
$ws.Range("M33").Value = -48110

# Row 35
$ws.Range("D35").Value = 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 74687
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0

# Row 36
$ws.Range("D36").Value = "-"
$ws.Range("E36").Value = "-"
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = "-"
$ws.Range("H36").Value = "-"
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = "-"
$ws.Range("L36").Value = "-"
$ws.Range("M36").Value = 0

# Row 37
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 427047
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 144944

# Row 38
$ws.Range("D38").Value = -45424
$ws.Range("E38").Value = -83467
$ws.Range("F38").Value = -62018
$ws.Range("G38").Value = -2432
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = -391452
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -159689

# Row 39
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 130000
$ws.Range("G39").Value = 200000
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 70000
$ws.Range("J39").Value = 150000
$ws.Range("K39").Value = 120000
$ws.Range("L39").Value = 150000
$ws.Range("M39").Value = 0

# Row 40
$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = -17887
$ws.Range("G40").Value = -62200
$ws.Range("H40").Value = -22329
$ws.Range("I40").Value = -153968
$ws.Range("J40").Value = -93890
$ws.Range("K40").Value = -89402
$ws.Range("L40").Value = -87827
$ws.Range("M40").Value = -163969

# Row 41
$ws.Range("D41").Value = 0
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = -1625
$ws.Range("G41").Value = -6524
$ws.Range("H41").Value = -10607
$ws.Range("I41").Value = -7418
$ws.Range("J41").Value = -6727
$ws.Range("K41").Value = -6207
$ws.Range("L41").Value = -9642
$ws.Range("M41").Value = -11024

# Row 42
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0

# Row 43
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0

# Row 44
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0

# Row 45
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0

# Row 46
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = 0

# Row 47
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0

# Row 48
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0

# Row 49
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0

# Row 50
$ws.Range("D50").Value = -102341
$ws.Range("E50").Value = -33345
$ws.Range("F50").Value = -170700
$ws.Range("G50").Value = -396683
$ws.Range("H50").Value = -302923
$ws.Range("I50").Value = -4850
$ws.Range("J50").Value = -112
$ws.Range("K50").Value = -552116
$ws.Range("L50").Value = -910713
$ws.Range("M50").Value = -443

# Row 51
$ws.Range("D51").Value = -147765
$ws.Range("E51").Value = -116812
$ws.Range("F51").Value = -122230
$ws.Range("G51").Value = -267839
$ws.Range("H51").Value = -335859
$ws.Range("I51").Value = -60641
$ws.Range("J51").Value = 49271
$ws.Range("K51").Value = -453038
$ws.Range("L51").Value = -858182
$ws.Range("M51").Value = -190181

# Row 52
$ws.Range("D52").Value = 25033
$ws.Range("E52").Value = 15338
$ws.Range("F52").Value = 152341
$ws.Range("G52").Value = -36650
$ws.Range("H52").Value = 21399
$ws.Range("I52").Value = -95535
$ws.Range("J52").Value = 249370
$ws.Range("K52").Value = -146188
$ws.Range("L52").Value = 68148
$ws.Range("M52").Value = -238291

# Row 53
$ws.Range("D53").Value = 36151
$ws.Range("E53").Value = 61184
$ws.Range("F53").Value = 76456
$ws.Range("G53").Value = 228797
$ws.Range("H53").Value = 192186
$ws.Range("I53").Value = 213585
$ws.Range("J53").Value = 118011
$ws.Range("K53").Value = 367381
$ws.Range("L53").Value = 221193
$ws.Range("M53").Value = 289341

# Row 54
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = -66
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 39
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = -39
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = 149

# Row 55
$ws.Range("D55").Value = 61184
$ws.Range("E55").Value = 76456
$ws.Range("F55").Value = 228797
$ws.Range("G55").Value = 192186
$ws.Range("H55").Value = 213585
$ws.Range("I55").Value = 118011
$ws.Range("J55").Value = 367381
$ws.Range("K55").Value = 221193
$ws.Range("L55").Value = 289341
$ws.Range("M55").Value = 51199

# Row 56
$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 0
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 14
$ws.Range("I56").Value = -14
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 313
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = 0

# Column width adjustments: the "year-end" quarter column (the widest one)
# shifts left by one position along with the data.
$ws.Columns.Item(4).ColumnWidth = 28.16666666666667   # D  width 29
$ws.Columns.Item(5).ColumnWidth = 30.16666666666667   # E  width 31
$ws.Columns.Item(6).ColumnWidth = 28.16666666666667   # F  width 29
$ws.Columns.Item(7).ColumnWidth = 28.16666666666667   # G  width 29
$ws.Columns.Item(8).ColumnWidth = 28.16666666666667   # H  width 29
$ws.Columns.Item(9).ColumnWidth = 30.16666666666667   # I  width 31
$ws.Columns.Item(10).ColumnWidth = 28.16666666666667  # J  width 29
$ws.Columns.Item(11).ColumnWidth = 28.16666666666667  # K  width 29
$ws.Columns.Item(12).ColumnWidth = 28.16666666666667  # L  width 29
$ws.Columns.Item(13).ColumnWidth = 30.16666666666667  # M  width 31
